# Workbook "Hortaliza, Vega Monumental Concepción - Perejil"
# The source data set gained one more weekly observation (commit: "Fruta / hortaliza, semanal").
# This shifts the existing data block (rows 144-207) down by two rows (146-209)
# and inserts a brand new pair of rows (144-145) at the top of the block holding
# the newest observation (Primera / Segunda quality) dated 2023-03-20 (serial 44992).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 144, pushing the
# existing rows 144:207 down to 146:209 (formatting/styles travel with them).
$ws.Rows.Item(144).Insert()
$ws.Rows.Item(144).Insert()

# New row 144 - "Primera" quality observation
$ws.Range("A144").Value = 11
$ws.Range("B144").Value = "Vega Monumental Concepción"
$ws.Range("C144").Value = "Bíobío"
$ws.Range("D144").Value = 44992
$ws.Range("E144").Value = 8
$ws.Range("F144").Value = 100112044
$ws.Range("G144").Value = "Perejil"
$ws.Range("H144").Value = "Sin especificar"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 200
$ws.Range("K144").Value = 700
$ws.Range("L144").Value = 800
$ws.Range("M144").Value = 750
$ws.Range("N144").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O144").Value = "Región de Ñuble"
$ws.Range("P144").Value = 750
$ws.Range("Q144").Value = 1
$ws.Range("R144").Value = "Hortaliza"

# New row 145 - "Segunda" quality observation
$ws.Range("A145").Value = 11
$ws.Range("B145").Value = "Vega Monumental Concepción"
$ws.Range("C145").Value = "Bíobío"
$ws.Range("D145").Value = 44992
$ws.Range("E145").Value = 8
$ws.Range("F145").Value = 100112044
$ws.Range("G145").Value = "Perejil"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Segunda"
$ws.Range("J145").Value = 100
$ws.Range("K145").Value = 600
$ws.Range("L145").Value = 600
$ws.Range("M145").Value = 600
$ws.Range("N145").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O145").Value = "Región de Ñuble"
$ws.Range("P145").Value = 600
$ws.Range("Q145").Value = 1
$ws.Range("R145").Value = "Hortaliza"
